# Updated cryptos list on Mon Jun  3 10:41:45 UTC 2024 with GitHub Actions
#
# Writes each changed Coin/Link/Price/Volume(1h) cell to match the refreshed
# coinranking.com snapshot. Price (column D) values that look like plain
# numbers are written with a leading apostrophe so Excel stores them as text
# (matching the original inlineStr cells: e.g. "1.00" must stay "1.00", not
# become the number 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($cellRef).Value = "'" + $val
    } else {
        $ws.Range($cellRef).Value = $val
    }
}

# --- Row 2 (Bitcoin) ---
Set-TextValue $ws "D2" "69.133.15"
Set-TextValue $ws "E2" "  +2.42%  "

# --- Row 3 (Ethereum) ---
Set-TextValue $ws "D3" "3.816.71"
Set-TextValue $ws "E3" "  +1.10%  "

# --- Row 4 (TetherUSD) ---
Set-TextValue $ws "E4" "  +0.12%  "

# --- Row 5 (BNB) ---
Set-TextValue $ws "D5" "630.47"
Set-TextValue $ws "E5" "  +5.46%  "

# --- Row 6 (Solana) ---
Set-TextValue $ws "D6" "165.46"
Set-TextValue $ws "E6" "  +0.73%  "

# --- Row 7 (LidoStakedEther) ---
Set-TextValue $ws "D7" "3.813.79"
Set-TextValue $ws "E7" "  +1.07%  "

# --- Row 8 (USDC) ---
Set-TextValue $ws "D8" "0.999"
Set-TextValue $ws "E8" "  -0.14%  "

# --- Row 9 (XRP) ---
Set-TextValue $ws "E9" "  +1.01%  "

# --- Row 10 (Dogecoin) ---
Set-TextValue $ws "D10" "0.163"
Set-TextValue $ws "E10" "  +3.07%  "

# --- Row 11 (Cardano) ---
Set-TextValue $ws "E11" "  +1.09%  "

# --- Row 12 (Toncoin) ---
Set-TextValue $ws "D12" "6.62"
Set-TextValue $ws "E12" "  +3.29%  "

# --- Row 13 (ShibaInu) ---
Set-TextValue $ws "E13" "  +1.19%  "

# --- Row 14 (Avalanche) ---
Set-TextValue $ws "D14" "36.07"
Set-TextValue $ws "E14" "  +1.74%  "

# --- Row 15 (WrappedliquidstakedEther2.0) ---
Set-TextValue $ws "D15" "4.457.49"
Set-TextValue $ws "E15" "  +1.12%  "

# --- Row 16 (WrappedEther) ---
Set-TextValue $ws "D16" "3.840.12"
Set-TextValue $ws "E16" "  +0.52%  "

# --- Row 17 (WrappedBTC) ---
Set-TextValue $ws "D17" "69.108.48"
Set-TextValue $ws "E17" "  +2.32%  "

# --- Row 18 (Chainlink) ---
Set-TextValue $ws "D18" "17.99"
Set-TextValue $ws "E18" "  -1.37%  "

# --- Row 19 (Polkadot) ---
Set-TextValue $ws "E19" "  +1.54%  "

# --- Row 21 (BitcoinCash) ---
Set-TextValue $ws "D21" "466.44"
Set-TextValue $ws "E21" "  +1.59%  "

# --- Row 22 (Uniswap) ---
Set-TextValue $ws "D22" "9.72"
Set-TextValue $ws "E22" "  +0.29%  "

# --- Row 23 (Polygon) ---
Set-TextValue $ws "D23" "0.709"
Set-TextValue $ws "E23" "  +2.22%  "

# --- Row 24 (PEPE) ---
Set-TextValue $ws "E24" "  +4.52%  "

# --- Row 25 (Litecoin) ---
Set-TextValue $ws "D25" "83.62"
Set-TextValue $ws "E25" "  +1.59%  "

# --- Row 26 (InternetComputer(DFINITY)) ---
Set-TextValue $ws "E26" "  +0.25%  "

# --- Row 27 (Fetch.AI) ---
Set-TextValue $ws "D27" "2.15"
Set-TextValue $ws "E27" "  +3.46%  "

# --- Row 28 (RenderToken) ---
Set-TextValue $ws "D28" "10.07"
Set-TextValue $ws "E28" "  +1.40%  "

# --- Row 29 (Dai) ---
Set-TextValue $ws "E29" "  +0.04%  "

# --- Row 30 (WrappedeETH) ---
Set-TextValue $ws "D30" "3.967.07"
Set-TextValue $ws "E30" "  +1.16%  "

# --- Row 31 (PancakeSwap) ---
Set-TextValue $ws "E31" "  +4.11%  "

# --- Row 32 (ImmutableX) ---
Set-TextValue $ws "D32" "2.22"
Set-TextValue $ws "E32" "  +1.88%  "

# --- Row 33 (NEARProtocol) ---
Set-TextValue $ws "E33" "  -1.79%  "

# --- Row 34 (EthereumClassic) ---
Set-TextValue $ws "D34" "29.23"
Set-TextValue $ws "E34" "  +0.67%  "

# --- Rows 35/36 swap: Binance-PegBSC-USD moves above Aptos, Aptos price refreshed ---
Set-TextValue $ws "B35" "Binance-PegBSC-USD"
Set-TextValue $ws "C35" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws "D35" "1.00"
Set-TextValue $ws "E35" "  +0.19%  "

Set-TextValue $ws "B36" "Aptos"
Set-TextValue $ws "C36" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D36" "9.08"
Set-TextValue $ws "E36" "  +1.52%  "

# --- Row 37 (Hedera) ---
Set-TextValue $ws "E37" "  +4.01%  "

# --- Row 38 (Kaspa) ---
Set-TextValue $ws "E38" "  +8.30%  "

# --- Row 39 (dogwifhat) ---
Set-TextValue $ws "D39" "3.44"
Set-TextValue $ws "E39" "  +6.56%  "

# --- Row 40 (Filecoin) ---
Set-TextValue $ws "E40" "  +3.22%  "

# --- Row 41 (Mantle) ---
Set-TextValue $ws "D41" "0.980"
Set-TextValue $ws "E41" "  -0.40%  "

# --- Row 42 (FirstDigitalUSD) ---
Set-TextValue $ws "E42" "  +0.11%  "

# --- Row 44 (Monero) ---
Set-TextValue $ws "D44" "158.00"
Set-TextValue $ws "E44" "  +4.29%  "

# --- Row 45 (ONDO) ---
Set-TextValue $ws "E45" "  +6.46%  "

# --- Row 46 (TheGraph) ---
Set-TextValue $ws "E46" "  +1.38%  "

# --- Row 47 (OKB) ---
Set-TextValue $ws "D47" "46.90"
Set-TextValue $ws "E47" "  -0.94%  "

# --- Row 48 (Arweave) ---
Set-TextValue $ws "D48" "42.80"
Set-TextValue $ws "E48" "  -1.53%  "

# --- Rows 49/50 swap: Stacks moves above Cosmos, Cosmos price refreshed ---
Set-TextValue $ws "B49" "Stacks"
Set-TextValue $ws "C49" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D49" "1.90"
Set-TextValue $ws "E49" "  +3.32%  "

Set-TextValue $ws "B50" "Cosmos"
Set-TextValue $ws "C50" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D50" "8.45"
Set-TextValue $ws "E50" "  +1.77%  "

# --- Row 51 (FLOKI) ---
Set-TextValue $ws "D51" "0.000284"
Set-TextValue $ws "E51" "  +15.38%  "
